# Mark the eight "messages.*" API methods that got implemented since the
# last update: they move from "In progress" to "Implemeted" on the
# Methods sheet (rows 42-49), picking up the same look (fill / number
# format / alignment) already used by the other "Implemeted" rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

for ($r = 42; $r -le 49; $r++) {
    $c = $ws.Cells.Item($r, 3)
    $c.Value = "Implemeted"
    $c.Style = "Акцент6"
    $c.NumberFormat = "@"
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4108
}

# messages.getByConversationMessageId (row 46) gets a side note.
$ws.Cells.Item(46, 4).Value = "?"

# Restore the scroll position / selection recorded the next time the
# workbook was saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("G57").Select()

$wb.Save()
